$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1 / tab 1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 56
$ws1.Range("F4").Value = 154
$ws1.Range("F5").Value = 353
$ws1.Range("F6").Value = 5117
$ws1.Range("F8").Value = 5289
$ws1.Range("F9").Value = 608
$ws1.Range("F10").Value = 1343

# Sheet "全部类型" (sheet4 / tab 4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 56
$ws4.Range("F4").Value = 154
$ws4.Range("F6").Value = 353
$ws4.Range("F7").Value = 5117
$ws4.Range("F9").Value = 5289
$ws4.Range("F10").Value = 608
$ws4.Range("F11").Value = 1343
